# Add 2022-Q3 data
# -----------------------------------------------------------------
# 1) Insert a brand-new worksheet named "2022-Q3" right before the
#    existing "2022-Q2" sheet (tab order: 总计, 2022-Q3, 2022-Q2, ...)
# -----------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$q2Ref = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Ref)
$newSheet.Name = "2022-Q3"

# Match page margins used by every other sheet in the workbook
# (0.75in / 0.75in / 1in / 1in / 0.5in / 0.5in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Re-fetch the template sheet by name (NOT the stale pre-insert
# reference) so formatting copies resolve against the right sheet.
$template = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row formatting (bold / centered / bordered) plus
# the index-column formatting from the template sheet so the new
# sheet's style matches every other quarter sheet.
for ($col = 2; $col -le 8; $col++) {
    $template.Cells.Item(1, $col).Copy($newSheet.Cells.Item(1, $col))
}
for ($row = 2; $row -le 8; $row++) {
    $template.Cells.Item(2, 1).Copy($newSheet.Cells.Item($row, 1))
}

# Header row text
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Columns B, D, E, F, G hold text-typed values (fund codes / figures
# stored as text, matching the source data) -- force text formatting
# on those ranges before writing so numeric-looking strings (e.g.
# "004856") are not auto-coerced into numbers.
$newSheet.Range("B2:B8").NumberFormat = "@"
$newSheet.Range("D2:G8").NumberFormat = "@"

$rows = @(
    @(0, "159745", "国泰中证全指建筑材料ETF",   "7.92", "99.14", "3.07", "0.2431", 9),
    @(1, "004856", "广发中证全指建筑材料指数A", "7.66", "93.74", "2.88", "0.2206", 9),
    @(2, "004857", "广发中证全指建筑材料指数C", "6.12", "93.74", "2.88", "0.1763", 9),
    @(3, "516750", "富国中证全指建筑材料ETF",   "0.82", "98.46", "3.07", "0.0252", 9),
    @(4, "159787", "易方达中证全指建筑材料ETF", "0.17", "94.24", "2.85", "0.0048", 9),
    @(5, "011015", "嘉合锦元回报混合A",         "0.78", "20.08", "0.22", "0.0017", 8),
    @(6, "011016", "嘉合锦元回报混合C",         "0.18", "20.08", "0.22", "0.0004", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# -----------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert the 2022-Q3 row at
#    the top of the data (row 2) and push the rest down, refreshing
#    the running index column (A).
# -----------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 7,  0.67),
    @(1, "2022-Q2", 4,  0.65),
    @(2, "2022-Q1", 6,  0.71),
    @(3, "2021-Q4", 21, 8.66),
    @(4, "2021-Q3", 14, 0.82),
    @(5, "2021-Q2", 36, 2.83),
    @(6, "2020-Q4", 2,  0.55)
)

# Make sure the newly introduced row 8 index cell carries the same
# formatting (bold / centered / bordered) as the rest of column A.
$total.Cells.Item(2, 1).Copy($total.Cells.Item(8, 1))

$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# -----------------------------------------------------------------
# 3) Leave the summary sheet as the active sheet/selection, matching
#    the workbook's original default view state.
# -----------------------------------------------------------------
$total.Activate()
$total.Range("A1").Select()
